$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear existing content (old layout: A1:I2) before rebuilding the full table
$ws.Cells.Clear()

# The numeric-looking columns (totalRuns, totalBalls, total4s, total6s, sr) must stay
# text, exactly as the source data has them (e.g. "162.50", "57.14"), so format the
# data portion of those columns as Text before writing the values into them.
$ws.Range("G2:K8").NumberFormat = "@"

# Header row (row 1)
$ws.Range("A1").Value = "venue"
$ws.Range("B1").Value = "date"
$ws.Range("C1").Value = "result"
$ws.Range("D1").Value = "ownTeam"
$ws.Range("E1").Value = "oppTeam"
$ws.Range("F1").Value = "batsman"
$ws.Range("G1").Value = "totalRuns"
$ws.Range("H1").Value = "totalBalls"
$ws.Range("I1").Value = "total4s"
$ws.Range("J1").Value = "total6s"
$ws.Range("K1").Value = "sr"

# Data rows (rows 2-8)
# Row 2
$ws.Range("A2").Value = " Sharjah"
$ws.Range("B2").Value = " September 27 2020"
$ws.Range("C2").Value = "Royals won by 4 wickets (with 3 balls remaining)"
$ws.Range("D2").Value = "Rajasthan Royals"
$ws.Range("E2").Value = "Kings XI Punjab"
$ws.Range("F2").Value = "Jos Buttler †"
$ws.Range("G2").Value = "4"
$ws.Range("H2").Value = "7"
$ws.Range("I2").Value = "0"
$ws.Range("J2").Value = "0"
$ws.Range("K2").Value = "57.14"

# Row 3
$ws.Range("A3").Value = " Dubai (DSC)"
$ws.Range("B3").Value = " September 30 2020"
$ws.Range("C3").Value = "KKR won by 37 runs"
$ws.Range("D3").Value = "Rajasthan Royals"
$ws.Range("E3").Value = "Kolkata Knight Riders"
$ws.Range("F3").Value = "Jos Buttler †"
$ws.Range("G3").Value = "21"
$ws.Range("H3").Value = "16"
$ws.Range("I3").Value = "1"
$ws.Range("J3").Value = "2"
$ws.Range("K3").Value = "131.25"

# Row 4
$ws.Range("A4").Value = " Abu Dhabi"
$ws.Range("B4").Value = " October 03 2020"
$ws.Range("C4").Value = "RCB won by 8 wickets (with 5 balls remaining)"
$ws.Range("D4").Value = "Rajasthan Royals"
$ws.Range("E4").Value = "Royal Challengers Bangalore"
$ws.Range("F4").Value = "Jos Buttler †"
$ws.Range("G4").Value = "22"
$ws.Range("H4").Value = "12"
$ws.Range("I4").Value = "3"
$ws.Range("J4").Value = "1"
$ws.Range("K4").Value = "183.33"

# Row 5
$ws.Range("A5").Value = " Dubai (DSC)"
$ws.Range("B5").Value = " October 14 2020"
$ws.Range("C5").Value = "Capitals won by 13 runs"
$ws.Range("D5").Value = "Rajasthan Royals"
$ws.Range("E5").Value = "Delhi Capitals"
$ws.Range("F5").Value = "Jos Buttler †"
$ws.Range("G5").Value = "22"
$ws.Range("H5").Value = "9"
$ws.Range("I5").Value = "3"
$ws.Range("J5").Value = "1"
$ws.Range("K5").Value = "244.44"

# Row 6
$ws.Range("A6").Value = " Sharjah"
$ws.Range("B6").Value = " October 09 2020"
$ws.Range("C6").Value = "Capitals won by 46 runs"
$ws.Range("D6").Value = "Rajasthan Royals"
$ws.Range("E6").Value = "Delhi Capitals"
$ws.Range("F6").Value = "Jos Buttler †"
$ws.Range("G6").Value = "13"
$ws.Range("H6").Value = "8"
$ws.Range("I6").Value = "2"
$ws.Range("J6").Value = "0"
$ws.Range("K6").Value = "162.50"

# Row 7
$ws.Range("A7").Value = " Abu Dhabi"
$ws.Range("B7").Value = " October 06 2020"
$ws.Range("C7").Value = "Mumbai won by 57 runs"
$ws.Range("D7").Value = "Rajasthan Royals"
$ws.Range("E7").Value = "Mumbai Indians"
$ws.Range("F7").Value = "Jos Buttler †"
$ws.Range("G7").Value = "70"
$ws.Range("H7").Value = "44"
$ws.Range("I7").Value = "4"
$ws.Range("J7").Value = "5"
$ws.Range("K7").Value = "159.09"

# Row 8
$ws.Range("A8").Value = " Dubai (DSC)"
$ws.Range("B8").Value = " October 11 2020"
$ws.Range("C8").Value = "Royals won by 5 wickets (with 1 ball remaining)"
$ws.Range("D8").Value = "Rajasthan Royals"
$ws.Range("E8").Value = "Sunrisers Hyderabad"
$ws.Range("F8").Value = "Jos Buttler †"
$ws.Range("G8").Value = "16"
$ws.Range("H8").Value = "13"
$ws.Range("I8").Value = "1"
$ws.Range("J8").Value = "1"
$ws.Range("K8").Value = "123.07"

